$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.145.85"
$ws.Range("E2").Value = "  +0.42%  "
$ws.Range("D3").Value = "2.282.54"
$ws.Range("E3").Value = "  -0.69%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'318.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.09%  "
$ws.Range("D6").Value = "'100.92"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.90%  "
$ws.Range("D7").Value = "'0.627"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "'0.603"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.85%  "
$ws.Range("D10").Value = "'38.86"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.63%  "
$ws.Range("D11").Value = "'0.0896"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.55%  "
$ws.Range("E12").Value = "  -2.12%  "
$ws.Range("E13").Value = "  -0.28%  "
$ws.Range("D14").Value = "'0.952"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.31%  "
$ws.Range("E15").Value = "  -1.89%  "
$ws.Range("D16").Value = "2.628.81"
$ws.Range("E16").Value = "  -0.54%  "
$ws.Range("D17").Value = "2.284.51"
$ws.Range("E17").Value = "  +0.08%  "
$ws.Range("D18").Value = "42.135.55"
$ws.Range("E18").Value = "  +0.53%  "
$ws.Range("D19").Value = "'7.30"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.17%  "
$ws.Range("E20").Value = "  -0.43%  "
$ws.Range("D21").Value = "'12.67"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +28.64%  "
$ws.Range("D22").Value = "'72.75"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.16%  "
$ws.Range("D23").Value = "'3.53"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.56%  "
$ws.Range("D24").Value = "'267.87"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.86%  "
$ws.Range("E25").Value = "  -4.50%  "
$ws.Range("E26").Value = "  -0.64%  "
$ws.Range("D28").Value = "'2.29"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.89%  "
$ws.Range("E29").Value = "  -1.79%  "
$ws.Range("D30").Value = "'37.56"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.66%  "
$ws.Range("D31").Value = "'163.82"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.43%  "
$ws.Range("D32").Value = "'6.03"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.66%  "
$ws.Range("E33").Value = "  -2.29%  "
$ws.Range("E35").Value = "  -3.27%  "
$ws.Range("D36").Value = "'2.46"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -15.55%  "
$ws.Range("E37").Value = "  -1.45%  "
$ws.Range("E38").Value = "  +0.54%  "
$ws.Range("E39").Value = "  -4.68%  "
$ws.Range("D40").Value = "'3.68"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.28%  "
$ws.Range("E41").Value = "  +2.27%  "
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("D43").Value = "'68.19"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.15%  "
$ws.Range("E44").Value = "  -1.62%  "
$ws.Range("D45").Value = "'91.41"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -9.44%  "
$ws.Range("E46").Value = "  +0.78%  "
$ws.Range("D47").Value = "'12.08"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.16%  "
$ws.Range("D48").Value = "'79.18"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.73%  "
$ws.Range("D49").Value = "'8.93"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.66%  "
$ws.Range("D50").Value = "1.608.53"
$ws.Range("E50").Value = "  +4.44%  "
$ws.Range("D51").Value = "'5.19"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.10%  "
